{"js": "const pairs = [\n  [\"19+28=\", \"65-1=\"],\n  [\"77-71=\", \"17+47=\"],\n  [\"64-20=\", \"47+1=\"],\n  [\"24+5=\", \"75+13=\"],\n  [\"60-14=\", \"63+6=\"],\n  [\"31+20=\", \"57-0=\"],\n  [\"37+51=\", \"44-31=\"],\n  [\"62-8=\", \"66+7=\"],\n  [\"23-18=\", \"98-37=\"],\n  [\"96-45=\", \"26+10=\"],\n  [\"2+87=\", \"0+59=\"],\n  [\"68+9=\", \"21+28=\"],\n  [\"53-38=\", \"51+41=\"],\n  [\"88-44=\", \"73+18=\"],\n  [\"65-33=\", \"40+14=\"],\n  [\"70-23=\", \"48+45=\"],\n  [\"99-87=\", \"39+18=\"],\n  [\"31+62=\", \"24+8=\"],\n  [\"87+2=\", \"37-35=\"],\n  [\"4+89=\", \"90-42=\"],\n  [\"31-0=\", \"64-40=\"],\n  [\"57-32=\", \"43-3=\"],\n  [\"15-6=\", \"3+30=\"],\n  [\"46-20=\", \"53+29=\"],\n  [\"50-34=\", \"63+14=\"],\n  [\"11+51=\", \"8+88=\"],\n  [\"48+6=\", \"81-68=\"],\n  [\"52+0=\", \"19+16=\"],\n  [\"62+20=\", \"55+4=\"],\n  [\"79-26=\", \"36+54=\"],\n  [\"78-56=\", \"24+25=\"],\n  [\"41+29=\", \"77-47=\"],\n  [\"36-10=\", \"31+1=\"],\n  [\"58+2=\", \"32+10=\"],\n  [\"22-5=\", \"79+6=\"],\n  [\"74-3=\", \"46+20=\"],\n  [\"92-12=\", \"53+40=\"],\n  [\"11+37=\", \"51+44=\"],\n  [\"72-8=\", \"65-17=\"],\n  [\"25+16=\", \"73-4=\"],\n  [\"72+25=\", \"52-26=\"],\n  [\"38-27=\", \"48+14=\"],\n  [\"14+31=\", \"19+42=\"],\n  [\"36-22=\", \"6+17=\"],\n  [\"55-36=\", \"37+35=\"],\n  [\"73-59=\", \"57-46=\"],\n  [\"7+53=\", \"45+54=\"],\n  [\"6+62=\", \"76-59=\"],\n  [\"12+65=\", \"78-20=\"],\n  [\"56-3=\", \"42+46=\"],\n  [\"22-21=\", \"19+51=\"],\n  [\"95-73=\", \"0+94=\"],\n  [\"71-37=\", \"19+66=\"],\n  [\"47+47=\", \"88-6=\"],\n  [\"68-54=\", \"92-51=\"],\n  [\"27-18=\", \"19+32=\"],\n  [\"54+29=\", \"56+38=\"],\n  [\"5+80=\", \"72-14=\"],\n  [\"90-50=\", \"2+97=\"],\n  [\"87-37=\", \"58-30=\"],\n  [\"97-47=\", \"45-39=\"],\n  [\"12+9=\", \"33+0=\"],\n  [\"89-19=\", \"78-46=\"],\n  [\"34+23=\", \"67+28=\"],\n  [\"57+13=\", \"21+60=\"],\n  [\"98-49=\", \"83-6=\"],\n  [\"89-16=\", \"22+39=\"],\n  [\"25-1=\", \"39+12=\"],\n  [\"10+84=\", \"69-23=\"],\n  [\"41+14=\", \"82-12=\"],\n  [\"35+1=\", \"64-27=\"],\n  [\"5+72=\", \"0+55=\"],\n  [\"60+2=\", \"93-58=\"],\n  [\"98-43=\", \"40+56=\"],\n  [\"86-73=\", \"12+68=\"],\n  [\"46+32=\", \"26+22=\"],\n  [\"84-19=\", \"56-6=\"],\n  [\"27+13=\", \"45+10=\"],\n  [\"13+18=\", \"45-22=\"],\n  [\"13+85=\", \"6-6=\"],\n  [\"33-16=\", \"37+62=\"],\n  [\"7+46=\", \"77-28=\"],\n  [\"45-7=\", \"76-24=\"],\n  [\"73-36=\", \"83-56=\"],\n  [\"67-63=\", \"0+31=\"],\n  [\"44-20=\", \"65-52=\"],\n  [\"22+17=\", \"37-16=\"],\n  [\"95-91=\", \"54-49=\"],\n  [\"7-6=\", \"60+34=\"],\n  [\"24+10=\", \"60-55=\"],\n  [\"16+79=\", \"36+0=\"],\n  [\"55+32=\", \"34-33=\"],\n  [\"84-9=\", \"46-21=\"],\n  [\"47+21=\", \"12+43=\"],\n  [\"14+55=\", \"29+30=\"],\n  [\"82-54=\", \"71-58=\"],\n  [\"96-17=\", \"28-9=\"],\n  [\"10+75=\", \"5+2=\"],\n  [\"98-26=\", \"17+30=\"],\n  [\"77-53=\", \"86-23=\"],\n];\n\nfor (const [oldText, newText] of pairs) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load('items');\n  await context.sync();\n  if (results.items.length === 0) {\n    throw new Error('No match found for: ' + oldText);\n  }\n  for (const item of results.items) {\n    item.insertText(newText, Word.InsertLocation.replace);\n  }\n}\nawait context.sync();", "ps1": "$d = $word.ActiveDocument\n$wdFindContinue = 1\n$wdReplaceAll = 2\n\n$pairs = @(\n  ,@(\"19+28=\", \"65-1=\")\n  ,@(\"77-71=\", \"17+47=\")\n  ,@(\"64-20=\", \"47+1=\")\n  ,@(\"24+5=\", \"75+13=\")\n  ,@(\"60-14=\", \"63+6=\")\n  ,@(\"31+20=\", \"57-0=\")\n  ,@(\"37+51=\", \"44-31=\")\n  ,@(\"62-8=\", \"66+7=\")\n  ,@(\"23-18=\", \"98-37=\")\n  ,@(\"96-45=\", \"26+10=\")\n  ,@(\"2+87=\", \"0+59=\")\n  ,@(\"68+9=\", \"21+28=\")\n  ,@(\"53-38=\", \"51+41=\")\n  ,@(\"88-44=\", \"73+18=\")\n  ,@(\"65-33=\", \"40+14=\")\n  ,@(\"70-23=\", \"48+45=\")\n  ,@(\"99-87=\", \"39+18=\")\n  ,@(\"31+62=\", \"24+8=\")\n  ,@(\"87+2=\", \"37-35=\")\n  ,@(\"4+89=\", \"90-42=\")\n  ,@(\"31-0=\", \"64-40=\")\n  ,@(\"57-32=\", \"43-3=\")\n  ,@(\"15-6=\", \"3+30=\")\n  ,@(\"46-20=\", \"53+29=\")\n  ,@(\"50-34=\", \"63+14=\")\n  ,@(\"11+51=\", \"8+88=\")\n  ,@(\"48+6=\", \"81-68=\")\n  ,@(\"52+0=\", \"19+16=\")\n  ,@(\"62+20=\", \"55+4=\")\n  ,@(\"79-26=\", \"36+54=\")\n  ,@(\"78-56=\", \"24+25=\")\n  ,@(\"41+29=\", \"77-47=\")\n  ,@(\"36-10=\", \"31+1=\")\n  ,@(\"58+2=\", \"32+10=\")\n  ,@(\"22-5=\", \"79+6=\")\n  ,@(\"74-3=\", \"46+20=\")\n  ,@(\"92-12=\", \"53+40=\")\n  ,@(\"11+37=\", \"51+44=\")\n  ,@(\"72-8=\", \"65-17=\")\n  ,@(\"25+16=\", \"73-4=\")\n  ,@(\"72+25=\", \"52-26=\")\n  ,@(\"38-27=\", \"48+14=\")\n  ,@(\"14+31=\", \"19+42=\")\n  ,@(\"36-22=\", \"6+17=\")\n  ,@(\"55-36=\", \"37+35=\")\n  ,@(\"73-59=\", \"57-46=\")\n  ,@(\"7+53=\", \"45+54=\")\n  ,@(\"6+62=\", \"76-59=\")\n  ,@(\"12+65=\", \"78-20=\")\n  ,@(\"56-3=\", \"42+46=\")\n  ,@(\"22-21=\", \"19+51=\")\n  ,@(\"95-73=\", \"0+94=\")\n  ,@(\"71-37=\", \"19+66=\")\n  ,@(\"47+47=\", \"88-6=\")\n  ,@(\"68-54=\", \"92-51=\")\n  ,@(\"27-18=\", \"19+32=\")\n  ,@(\"54+29=\", \"56+38=\")\n  ,@(\"5+80=\", \"72-14=\")\n  ,@(\"90-50=\", \"2+97=\")\n  ,@(\"87-37=\", \"58-30=\")\n  ,@(\"97-47=\", \"45-39=\")\n  ,@(\"12+9=\", \"33+0=\")\n  ,@(\"89-19=\", \"78-46=\")\n  ,@(\"34+23=\", \"67+28=\")\n  ,@(\"57+13=\", \"21+60=\")\n  ,@(\"98-49=\", \"83-6=\")\n  ,@(\"89-16=\", \"22+39=\")\n  ,@(\"25-1=\", \"39+12=\")\n  ,@(\"10+84=\", \"69-23=\")\n  ,@(\"41+14=\", \"82-12=\")\n  ,@(\"35+1=\", \"64-27=\")\n  ,@(\"5+72=\", \"0+55=\")\n  ,@(\"60+2=\", \"93-58=\")\n  ,@(\"98-43=\", \"40+56=\")\n  ,@(\"86-73=\", \"12+68=\")\n  ,@(\"46+32=\", \"26+22=\")\n  ,@(\"84-19=\", \"56-6=\")\n  ,@(\"27+13=\", \"45+10=\")\n  ,@(\"13+18=\", \"45-22=\")\n  ,@(\"13+85=\", \"6-6=\")\n  ,@(\"33-16=\", \"37+62=\")\n  ,@(\"7+46=\", \"77-28=\")\n  ,@(\"45-7=\", \"76-24=\")\n  ,@(\"73-36=\", \"83-56=\")\n  ,@(\"67-63=\", \"0+31=\")\n  ,@(\"44-20=\", \"65-52=\")\n  ,@(\"22+17=\", \"37-16=\")\n  ,@(\"95-91=\", \"54-49=\")\n  ,@(\"7-6=\", \"60+34=\")\n  ,@(\"24+10=\", \"60-55=\")\n  ,@(\"16+79=\", \"36+0=\")\n  ,@(\"55+32=\", \"34-33=\")\n  ,@(\"84-9=\", \"46-21=\")\n  ,@(\"47+21=\", \"12+43=\")\n  ,@(\"14+55=\", \"29+30=\")\n  ,@(\"82-54=\", \"71-58=\")\n  ,@(\"96-17=\", \"28-9=\")\n  ,@(\"10+75=\", \"5+2=\")\n  ,@(\"98-26=\", \"17+30=\")\n  ,@(\"77-53=\", \"86-23=\")\n)\n\nforeach ($pair in $pairs) {\n  $oldText = $pair[0]\n  $newText = $pair[1]\n  $find = $d.Content.Find\n  $find.ClearFormatting()\n  $find.Replacement.ClearFormatting()\n  $find.Text = $oldText\n  $find.Replacement.Text = $newText\n  $find.Execute($oldText, $false, $false, $false, $false, $false, $true, $wdFindContinue, $false, $newText, $wdReplaceAll) | Out-Null\n}"}
